$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.649.22"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'1.613.12"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.63%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.994"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.58%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'212.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.15%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.520"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.34%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.993"
$ws.Range("D7").Style = "Normal"
$ws.Range("E8").Value = "'  +9.48%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +3.01%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.0608"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.61%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0908"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.13%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.847.71"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.81%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.612.74"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.51%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +6.96%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.86"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +4.84%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'29.668.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.56%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'8.96"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +17.44%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'64.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.36%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'241.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.37%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +3.00%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  -0.43%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +2.77%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'9.67"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +5.73%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +1.40%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'156.08"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +1.14%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'15.62"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +2.25%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +2.12%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +3.65%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -0.49%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.0487"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +3.23%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("B31").Value = "'PancakeSwap"
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value = "'1.08"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +2.31%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("B32").Value = "'Filecoin"
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = "'3.31"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +3.09%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'3.20"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +3.37%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.433.54"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.33%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.62"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +6.74%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +1.65%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +2.17%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'2.28"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -1.05%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +3.16%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.557"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +4.09%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.0499"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +3.09%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.829"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +4.51%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +0.53%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'69.94"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +6.65%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'53.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.12%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.993"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.55%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.999"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +18.41%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +3.32%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.755.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'87.92"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.59%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'2.11"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.04%  "
$ws.Range("E51").Style = "Normal"
